# Fix bugs accessing GBIF, improve output
#
# Sheet1 (the active sheet) had the GAVIA dataset description in row 3 and the
# GloNAF dataset description in row 5. The two rows' contents are swapped so
# that row 3 now documents GloNAF and row 5 documents GAVIA. The GloNAF row
# additionally uses columns G ("author") and K ("status"), which the GAVIA row
# does not use, so those two cells are cleared on row 5 once its content
# reverts to the (shorter) GAVIA description. Sheet2 already lists GloNAF in
# row 5 and GAVIA in row 3, so it is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Capture current ("before") values of the two rows that swap places ---

# Row 3 currently holds the GAVIA entry
$row3_A = $ws.Range("A3").Value2
$row3_B = $ws.Range("B3").Value2
$row3_C = $ws.Range("C3").Value2
$row3_E = $ws.Range("E3").Value2
$row3_I = $ws.Range("I3").Value2

# Row 5 currently holds the GloNAF entry
$row5_A = $ws.Range("A5").Value2
$row5_B = $ws.Range("B5").Value2
$row5_C = $ws.Range("C5").Value2
$row5_E = $ws.Range("E5").Value2
$row5_G = $ws.Range("G5").Value2
$row5_I = $ws.Range("I5").Value2
$row5_K = $ws.Range("K5").Value2

# --- Write the GloNAF entry (old row 5) into row 3 ---
$ws.Range("A3").Value2 = $row5_A
$ws.Range("B3").Value2 = $row5_B
$ws.Range("C3").Value2 = $row5_C
$ws.Range("E3").Value2 = $row5_E
$ws.Range("G3").Value2 = $row5_G
$ws.Range("I3").Value2 = $row5_I
$ws.Range("K3").Value2 = $row5_K

# --- Write the GAVIA entry (old row 3) into row 5 ---
$ws.Range("A5").Value2 = $row3_A
$ws.Range("B5").Value2 = $row3_B
$ws.Range("C5").Value2 = $row3_C
$ws.Range("E5").Value2 = $row3_E
$ws.Range("I5").Value2 = $row3_I

# Row 5 is now the (shorter) GAVIA entry, which has nothing in columns G/K
$ws.Range("G5").ClearContents()
$ws.Range("K5").ClearContents()

# Update the active selection on the active sheet to A5
$ws.Range("A5").Select()
